$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Update the environment URL cells (A2/B2) to the "i-" prefixed preprod host
$ws.Range("A2").Value = "i-preproducciongestion.segurossura.com.ar"
$ws.Range("B2").Value = "https://i-preproducciongestion.segurossura.com.ar/pc/PolicyCenter.do"

# Update the account number used for the regression test data
$ws.Range("E2").Value = 7068873718

# Move the active selection to Q2 (was Q5)
$ws.Range("Q2").Select()
